$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "29.098.81"
Set-TextValue "E2" "  +1.14%  "
Set-TextValue "D3" "1.924.26"
Set-TextValue "E3" "  +2.03%  "
Set-TextValue "E4" "  +0.41%  "
Set-TextValue "D5" "326.15"
Set-TextValue "E5" "  +1.12%  "
Set-TextValue "E6" "  +0.46%  "
Set-TextValue "D7" "0.4611"
Set-TextValue "E7" "  +1.05%  "
Set-TextValue "D8" "0.3823"
Set-TextValue "E8" "  +0.62%  "
Set-TextValue "D9" "0.07769"
Set-TextValue "E9" "  +0.78%  "
Set-TextValue "D10" "0.9790"
Set-TextValue "E10" "  +1.92%  "
Set-TextValue "D11" "22.66"
Set-TextValue "E11" "  +3.15%  "
Set-TextValue "D12" "1.930.03"
Set-TextValue "E12" "  +2.13%  "
Set-TextValue "D13" "5.698"
Set-TextValue "E13" "  +0.79%  "
Set-TextValue "D14" "6.971"
Set-TextValue "E14" "  +0.42%  "
Set-TextValue "D15" "0.07046"
Set-TextValue "E15" "  +0.83%  "
Set-TextValue "E16" "  +0.43%  "
Set-TextValue "D17" "84.38"
Set-TextValue "E17" "  +1.43%  "
Set-TextValue "D18" "0.000009531"
Set-TextValue "E18" "  +0.68%  "
Set-TextValue "D19" "16.73"
Set-TextValue "E19" "  +1.08%  "
Set-TextValue "E20" "  +0.47%  "
Set-TextValue "D21" "29.108.53"
Set-TextValue "E21" "  +1.37%  "
Set-TextValue "D22" "5.343"
Set-TextValue "E22" "  +0.54%  "
Set-TextValue "D23" "10.97"
Set-TextValue "E23" "  +1.06%  "
Set-TextValue "E24" "  -0.27%  "
Set-TextValue "D25" "158.15"
Set-TextValue "E25" "  +1.74%  "
Set-TextValue "E26" "  +1.00%  "
Set-TextValue "D27" "5.661"
Set-TextValue "E27" "  +1.25%  "
Set-TextValue "D28" "118.08"
Set-TextValue "E28" "  +1.09%  "
Set-TextValue "D29" "1.837"
Set-TextValue "E29" "  +2.11%  "
Set-TextValue "E30" "  +1.29%  "
Set-TextValue "D31" "0.8579"
Set-TextValue "E31" "  +1.79%  "
Set-TextValue "D32" "5.114"
Set-TextValue "E32" "  +1.11%  "
Set-TextValue "D33" "1.242"
Set-TextValue "E33" "  +0.55%  "
Set-TextValue "D34" "3.023"
Set-TextValue "E34" "  +1.06%  "
Set-TextValue "E35" "  +1.74%  "
Set-TextValue "D36" "0.05691"
Set-TextValue "E36" "  +0.58%  "
Set-TextValue "D37" "3.189"
Set-TextValue "E37" "  +18.26%  "
Set-TextValue "E38" "  +0.48%  "
Set-TextValue "E39" "  +1.02%  "
Set-TextValue "D40" "7.500"
Set-TextValue "E40" "  +1.27%  "
Set-TextValue "D41" "0.5514"
Set-TextValue "E41" "  +0.70%  "
Set-TextValue "D42" "0.1757"
Set-TextValue "E42" "  +0.60%  "
Set-TextValue "D43" "9.332"
Set-TextValue "E43" "  +2.28%  "
Set-TextValue "D44" "2.191"
Set-TextValue "E44" "  +6.50%  "
Set-TextValue "D45" "0.000002730"
Set-TextValue "E45" "  -8.38%  "
Set-TextValue "D46" "0.5207"
Set-TextValue "E46" "  +1.27%  "
Set-TextValue "B47" "EnergySwap"
Set-TextValue "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "11.27"
Set-TextValue "E47" "  +0.14%  "
Set-TextValue "B48" "Cronos"
Set-TextValue "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.06928"
Set-TextValue "E48" "  +1.97%  "
Set-TextValue "D49" "110.43"
Set-TextValue "E49" "  -0.92%  "
Set-TextValue "D50" "1.767"
Set-TextValue "E50" "  -0.13%  "
Set-TextValue "E51" "  +0.51%  "
